$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 290.08334
$ws.Range("I53").Value = 186.33333
$ws.Range("J53").Value = 393.83334
$ws.Range("K53").Value = 186.33333
$ws.Range("L53").Value = 393.83334
$ws.Range("M53").Value = 450.66667
$ws.Range("N53").Value = -1667.83334

$ws.Range("H76").Value = 5623
$ws.Range("I76").Value = 5623
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 5623
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -5308
$ws.Range("N76").Value = $null

$ws.Range("H79").Value = 5623
$ws.Range("I79").Value = 5623
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 5623
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -4531
$ws.Range("N79").Value = $null

$ws.Range("H82").Value = 484
$ws.Range("I82").Value = 484
$ws.Range("K82").Value = 1452
$ws.Range("M82").Value = -1046

$ws.Range("H85").Value = 484
$ws.Range("I85").Value = 484
$ws.Range("K85").Value = 1452
$ws.Range("M85").Value = -48

$ws.Range("H101").Value = 16667274
$ws.Range("J101").Value = 774.25
$ws.Range("L101").Value = 2322.75
$ws.Range("N101").Value = -5566.75

$ws.Range("H107").Value = 321.16666
$ws.Range("I107").Value = 321.16666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 321.16666
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1598.83334
$ws.Range("N107").Value = $null

$ws.Range("H132").Value = 4333.3335
$ws.Range("I132").Value = 4333.3335
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13000.0005
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10470.0005
$ws.Range("N132").Value = $null

$ws.Range("H140").Value = 55000
$ws.Range("J140").Value = 55000
$ws.Range("L140").Value = 55000
$ws.Range("N140").Value = -65360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1260.7142
$ws.Range("I2").Value = 365.4
$ws.Range("J2").Value = 3499
$ws.Range("K2").Value = 365.4
$ws.Range("L2").Value = 3499
$ws.Range("M2").Value = -252.4
$ws.Range("N2").Value = -3725

$ws.Range("H32").Value = 3503705
$ws.Range("I32").Value = 5838450
$ws.Range("K32").Value = 5838450
$ws.Range("M32").Value = -5838163

$ws.Range("H74").Value = 552.2
$ws.Range("I74").Value = 552.2
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 552.2
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 321.8
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 552.2
$ws.Range("I77").Value = 552.2
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 2761
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 1607
$ws.Range("N77").Value = $null

$ws.Range("H116").Value = 1260.7142
$ws.Range("I116").Value = 365.4
$ws.Range("J116").Value = 3499
$ws.Range("K116").Value = 365.4
$ws.Range("L116").Value = 3499
$ws.Range("M116").Value = 1928.6
$ws.Range("N116").Value = -8087

$ws.Range("H122").Value = 2180.7693
$ws.Range("I122").Value = 1704.5454
$ws.Range("K122").Value = 5113.6362
$ws.Range("M122").Value = -2663.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1260.7142
$ws.Range("I3").Value = 365.4
$ws.Range("J3").Value = 3499
$ws.Range("K3").Value = 365.4
$ws.Range("L3").Value = 3499
$ws.Range("M3").Value = -251.4
$ws.Range("N3").Value = -3727

$ws.Range("H134").Value = 2425.5
$ws.Range("I134").Value = 2425.5
$ws.Range("K134").Value = 7276.5
$ws.Range("M134").Value = -4741.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 81.15385000000001
$ws.Range("I7").Value = 52.42857
$ws.Range("J7").Value = 114.666664
$ws.Range("K7").Value = 52.42857
$ws.Range("L7").Value = 114.666664
$ws.Range("M7").Value = 60.57143
$ws.Range("N7").Value = -340.666664

$ws.Range("H31").Value = 727.75
$ws.Range("J31").Value = 505.5
$ws.Range("L31").Value = 505.5
$ws.Range("N31").Value = -1095.5

$ws.Range("H34").Value = 727.75
$ws.Range("J34").Value = 505.5
$ws.Range("L34").Value = 505.5
$ws.Range("N34").Value = -909.5

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H134").Value = 3033.3333
$ws.Range("I134").Value = 3033.3333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9099.999899999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6564.999899999999
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 162.72728
$ws.Range("I2").Value = 58.25
$ws.Range("J2").Value = 222.42857
$ws.Range("K2").Value = 349.5
$ws.Range("L2").Value = 1334.57142
$ws.Range("M2").Value = -236.5
$ws.Range("N2").Value = -1560.57142

$ws.Range("H11").Value = 116.111115
$ws.Range("I11").Value = 112.5
$ws.Range("K11").Value = 337.5
$ws.Range("M11").Value = -197.5

$ws.Range("H102").Value = 20000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 20000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 60000
$ws.Range("M102").Value = $null
$ws.Range("N102").Value = -64868

$ws.Range("H128").Value = 598339
$ws.Range("I128").Value = 598339
$ws.Range("K128").Value = 1795017
$ws.Range("M128").Value = -1790037

$ws.Range("H131").Value = 2349.1875
$ws.Range("I131").Value = 1172.75
$ws.Range("J131").Value = 2741.3333
$ws.Range("K131").Value = 3518.25
$ws.Range("L131").Value = 8223.999899999999
$ws.Range("M131").Value = 1521.75
$ws.Range("N131").Value = -18303.9999

$ws.Range("H133").Value = 3649.5
$ws.Range("I133").Value = 3649.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 10948.5
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -5888.5
$ws.Range("N133").Value = $null

$ws.Range("H134").Value = 1933
$ws.Range("I134").Value = 1933
$ws.Range("K134").Value = 5799
$ws.Range("M134").Value = -729

$ws.Range("H137").Value = 1920.091
$ws.Range("I137").Value = 1148.25
$ws.Range("J137").Value = 2361.1428
$ws.Range("K137").Value = 3444.75
$ws.Range("L137").Value = 7083.428400000001
$ws.Range("M137").Value = 1655.25
$ws.Range("N137").Value = -17283.4284

$ws.Range("H139").Value = 1270.2858
$ws.Range("I139").Value = 1198.5
$ws.Range("J139").Value = 1299
$ws.Range("K139").Value = 3595.5
$ws.Range("L139").Value = 3897
$ws.Range("M139").Value = 1544.5
$ws.Range("N139").Value = -14177

$ws.Range("H140").Value = 13818.615
$ws.Range("I140").Value = 3549
$ws.Range("K140").Value = 10647
$ws.Range("M140").Value = -5467

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 86.52
$ws.Range("I2").Value = 96.545456
$ws.Range("J2").Value = 13
$ws.Range("K2").Value = 96.545456
$ws.Range("L2").Value = 13
$ws.Range("M2").Value = 16.454544
$ws.Range("N2").Value = -239

$ws.Range("H34").Value = 43500
$ws.Range("I34").Value = 37000
$ws.Range("K34").Value = 37000
$ws.Range("M34").Value = -36732

$ws.Range("H74").Value = 49999.75
$ws.Range("J74").Value = 49999.75
$ws.Range("L74").Value = 49999.75
$ws.Range("N74").Value = -51871.75

$ws.Range("H76").Value = 43500
$ws.Range("I76").Value = 37000
$ws.Range("K76").Value = 37000
$ws.Range("M76").Value = -36685

$ws.Range("H77").Value = 49999.75
$ws.Range("J77").Value = 49999.75
$ws.Range("L77").Value = 149999.25
$ws.Range("N77").Value = -159359.25

$ws.Range("H79").Value = 43500
$ws.Range("I79").Value = 37000
$ws.Range("K79").Value = 37000
$ws.Range("M79").Value = -35908

$ws.Range("H113").Value = 610.1111

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 387.4
$ws.Range("I16").Value = 387.4
$ws.Range("K16").Value = 387.4
$ws.Range("M16").Value = -217.4

$ws.Range("H132").Value = 23096.6
$ws.Range("J132").Value = 28249.5
$ws.Range("L132").Value = 84748.5
$ws.Range("N132").Value = -89808.5

$ws.Range("H136").Value = 2984.4285
$ws.Range("I136").Value = 999.5
$ws.Range("K136").Value = 2998.5
$ws.Range("M136").Value = -448.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19979
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").Value = $null

$ws.Range("H75").Value = 74801.25
$ws.Range("J75").Value = 75399
$ws.Range("L75").Value = 75399
$ws.Range("N75").Value = -77271

$ws.Range("H78").Value = 74801.25
$ws.Range("J78").Value = 75399
$ws.Range("L78").Value = 226197
$ws.Range("N78").Value = -235557

$ws.Range("H113").Value = 244.71428
$ws.Range("I113").Value = 242.8
$ws.Range("J113").Value = 249.5
$ws.Range("K113").Value = 728.4000000000001
$ws.Range("L113").Value = 748.5
$ws.Range("N113").Value = -5088.5
$ws.Range("M113").Value = 1441.6

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550

$ws.Range("H131").Value = 14899.5
$ws.Range("J131").Value = 14899.5
$ws.Range("L131").Value = 14899.5
$ws.Range("N131").Value = -24979.5

$ws.Range("H132").Value = 4099.75
$ws.Range("I132").Value = 4799.6665
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 14398.9995
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -11868.9995
$ws.Range("N132").Value = -11060

$ws.Range("H136").Value = 3535.111
$ws.Range("I136").Value = 3802.4285
$ws.Range("J136").Value = 2599.5
$ws.Range("K136").Value = 11407.2855
$ws.Range("L136").Value = 7798.5
$ws.Range("M136").Value = -8857.2855
$ws.Range("N136").Value = -12898.5
